# feat: add 2022-Q3 data
#
# 1) Insert a brand-new "2022-Q3" sheet right after the "总计" summary sheet
#    (i.e. before "2021-Q4"), populated with the new quarter's holdings.
# 2) Update the "总计" sheet: the newest quarter (2022-Q3) becomes the new
#    row 2, and the previously-existing rows shift down to make room.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet    = $wb.Worksheets.Item("2021-Q4")

# --- 1) Create the new "2022-Q3" worksheet, positioned between 总计 and 2021-Q4 ---
$q3_2022 = $wb.Worksheets.Add($null, $totalSheet)
$q3_2022.Name = "2022-Q3"

$q3_2022.Range("B1").Value = "基金代码"
$q3_2022.Range("C1").Value = "基金名称"
$q3_2022.Range("D1").Value = "基金规模"
$q3_2022.Range("E1").Value = "股票总仓位"
$q3_2022.Range("F1").Value = "仓位占比"
$q3_2022.Range("G1").Value = "持有市值(亿元)"
$q3_2022.Range("H1").Value = "仓位排名"
$q3_2022.Range("B1:H1").Style = $q4Sheet.Range("B1:H1").Style

$q3_2022.Range("A2").Value = 0
$q3_2022.Range("B2").Value = "630010"
$q3_2022.Range("C2").Value = "华商价值精选混合"
$q3_2022.Range("D2").Value = "4.30"
$q3_2022.Range("E2").Value = "81.81"
$q3_2022.Range("F2").Value = "3.40"
$q3_2022.Range("G2").Value = "0.1462"
$q3_2022.Range("H2").Value = 6
$q3_2022.Range("A2").Style = $q4Sheet.Range("A2").Style

$q3_2022.Range("A3").Value = 1
$q3_2022.Range("B3").Value = "630006"
$q3_2022.Range("C3").Value = "华商产业升级混合"
$q3_2022.Range("D3").Value = "0.85"
$q3_2022.Range("E3").Value = "81.97"
$q3_2022.Range("F3").Value = "3.40"
$q3_2022.Range("G3").Value = "0.0289"
$q3_2022.Range("H3").Value = 6
$q3_2022.Range("A3").Style = $q4Sheet.Range("A2").Style

# --- 2) Update the "总计" (totals) sheet: shift rows down, insert new top row ---
$totalSheet.Range("A4").Value = $totalSheet.Range("A3").Value
$totalSheet.Range("B4").Value = $totalSheet.Range("B3").Value
$totalSheet.Range("C4").Value = $totalSheet.Range("C3").Value
$totalSheet.Range("D4").Value = $totalSheet.Range("D3").Value
$totalSheet.Range("A4").Style = $totalSheet.Range("A3").Style

$totalSheet.Range("A3").Value = $totalSheet.Range("A2").Value
$totalSheet.Range("B3").Value = $totalSheet.Range("B2").Value
$totalSheet.Range("C3").Value = $totalSheet.Range("C2").Value
$totalSheet.Range("D3").Value = $totalSheet.Range("D2").Value

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.03

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.18
